$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "Av. Salaverry 300"
$ws.Range("E6").Value = "Av. Salaverry 300, Candarave, Tacna, Tacna, Perú"
$ws.Range("F6").Value = -12.097018
$ws.Range("G6").Value = -77.055646
$ws.Range("H6").Value = "Oficina Centro"
